$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Overview sheet: handoff status text changed ("Ready for handoff" ->
#    "Handed back: in sync with en-US") for both language columns / rows.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# 2) zh-cn sheet: handback target/file columns + handback datetime filled in.
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b0ab24d9f82715506386ecce96dba969b61da611/e2e/6e5d81e4-f9d4-41fc-98c0-d323f2967b3e.md", $null, $null, "6e5d81e4-f9d4-41fc-98c0-d323f2967b3e.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b0ab24d9f82715506386ecce96dba969b61da611/e2e/6e5d81e4-f9d4-41fc-98c0-d323f2967b3e.md", $null, $null, "6e5d81e4-f9d4-41fc-98c0-d323f2967b3e.md") | Out-Null

$wsZh.Range("J2").Value = "6e5d81e4-f9d4-41fc-98c0-d323f2967b3e.5ab659c663c279dccc10836035f86cb1f1084335.zh-cn.xlf"
$wsZh.Range("J3").Value = "6e5d81e4-f9d4-41fc-98c0-d323f2967b3e.5ab659c663c279dccc10836035f86cb1f1084335.zh-cn.xlf"

$wsZh.Range("K2").Value = "2016-11-08 23:18:29"
$wsZh.Range("K3").Value = "2016-11-08 23:18:29"

# ---------------------------------------------------------------------------
# 3) de-de sheet: same pattern, but with the de-de target file + its own
#    handback datetime.
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b0ab24d9f82715506386ecce96dba969b61da611/e2e/6e5d81e4-f9d4-41fc-98c0-d323f2967b3e.md", $null, $null, "6e5d81e4-f9d4-41fc-98c0-d323f2967b3e.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b0ab24d9f82715506386ecce96dba969b61da611/e2e/6e5d81e4-f9d4-41fc-98c0-d323f2967b3e.md", $null, $null, "6e5d81e4-f9d4-41fc-98c0-d323f2967b3e.md") | Out-Null

$wsDe.Range("J2").Value = "6e5d81e4-f9d4-41fc-98c0-d323f2967b3e.5ab659c663c279dccc10836035f86cb1f1084335.de-de.xlf"
$wsDe.Range("J3").Value = "6e5d81e4-f9d4-41fc-98c0-d323f2967b3e.5ab659c663c279dccc10836035f86cb1f1084335.de-de.xlf"

$wsDe.Range("K2").Value = "2016-11-08 23:18:48"
$wsDe.Range("K3").Value = "2016-11-08 23:18:48"

# ---------------------------------------------------------------------------
# 4) Column widths: let Excel re-fit the columns that now hold longer text
#    (status column, and the newly populated target/handback-file columns).
# ---------------------------------------------------------------------------
$wsOverview.Range("E:F").Columns.AutoFit() | Out-Null
$wsZh.Range("C:C").Columns.AutoFit() | Out-Null
$wsZh.Range("I:J").Columns.AutoFit() | Out-Null
$wsDe.Range("C:C").Columns.AutoFit() | Out-Null
$wsDe.Range("I:J").Columns.AutoFit() | Out-Null

Write-Host "Generated handback report"
